# Edit script for 纱.xlsx
#
# What changed (per the target diff):
#   1. For each year's 4-quarter block, the "B" and "C" labelled rows had
#      their entire row content (columns A:E) swapped with one another.
#      The swapped row *pairs* (1-based worksheet rows) are:
#        (3,4) (7,8) (11,12) (15,16) (19,20) (23,24) (27,28) (31,32)
#        (35,36) (39,40) (43,44) (47,48) (51,52) (55,56) (59,60) (63,64) (67,68)
#   2. Columns F ("纱产销率") and G ("纱销售量") were deleted entirely
#      (including their header cells), shrinking the used range from
#      A1:G69 to A1:E69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(3,4), @(7,8), @(11,12), @(15,16), @(19,20), @(23,24), @(27,28),
    @(31,32), @(35,36), @(39,40), @(43,44), @(47,48), @(51,52),
    @(55,56), @(59,60), @(63,64), @(67,68)
)

$cols = @("A","B","C","D","E")

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $ref1 = $col + $r1
        $ref2 = $col + $r2

        $v1 = $ws.Range($ref1).Value2
        $v2 = $ws.Range($ref2).Value2

        # Skip cells that are empty on both sides so we don't disturb an
        # (intentionally present but blank) cell's stored representation.
        if (($v1 -eq "") -and ($v2 -eq "")) {
            continue
        }

        $ws.Range($ref1).Value = $v2
        $ws.Range($ref2).Value = $v1
    }
}

# Remove the now-unused columns F (纱产销率) and G (纱销售量) entirely.
$ws.Range("F1:G1").EntireColumn.Delete()
